# Apply the edits described by the commit:
# "added the code to compare the confirmation message text"
#
# - Adds 3 new rows (5, 6, 7) with new shared-string values, reusing the
#   same cell styles used by the existing header/label rows.
# - Widens column A slightly (separating it from columns B:C which keep
#   their original width).
# - Updates the active selection to O8 (as last left by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------

# Row 5 mirrors the style of row 1 (A1 = label style with border,
# B1/C1 = matching fill without border) but is left with empty B/C cells.
$ws.Range("A1:C1").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5").Value = "verifyProductAddedInCartConfirmationMessage"

# Row 6, column A mirrors the style used by row 2 (A2/B2).
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").Value = "Confirmation_Message"

# Row 7 uses the default (no explicit) style, like row 3.
$ws.Range("A7").Value = "Added to Cart"

# --- Column widths ---------------------------------------------------------
# Column A becomes its own width, separate from columns B:C which keep
# their existing shared width.
$ws.Columns.Item(1).ColumnWidth = 21.307291666666668

# --- Selection ---------------------------------------------------------
$ws.Range("O8").Select() | Out-Null
